$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = '顏之義(c_personid=32425)的籍貫是？'
$ws.Range('B2').Value = '長安'
$ws.Range('A3').Value = '朱廷桂(c_personid=62740的籍貫是否為東都指揮使司？'
$ws.Range('B3').Value = '否'
$ws.Range('A4').Value = '袁國鳳(c_personid=655539的籍貫是否為為州？'
$ws.Range('B4').Value = '否'
$ws.Range('A5').Value = '王有壬(c_personid=34586)的入仕方式是？'
$ws.Range('B5').Value = '恩蔭'
$ws.Range('A6').Value = '袁衡(c_personid=19396)的入仕方式是？'
$ws.Range('B6').Value = '進士'
$ws.Range('A7').Value = '李虛舟(c_personid=3505)的入仕方式是？'
$ws.Range('B7').Value = '恩蔭'
$ws.Range('A8').Value = '胡景定(c_personid=365873)的入仕方式是？'
$ws.Range('B8').Value = '進士'
$ws.Range('A9').Value = '衝有雅(c_personid=655421)的入仕方式是？'
$ws.Range('B9').Value = '舉人'
$ws.Range('A10').Value = '沈邦本(c_personid=632676的籍貫是否為寧州？'
$ws.Range('B10').Value = '否'
$ws.Range('A11').Value = '陳同熙(c_personid=669471)的入仕方式是？'
$ws.Range('B11').Value = '舉人'
$ws.Range('A12').Value = '陳大典(c_personid=371817)的入仕方式是？'
$ws.Range('B12').Value = '進士'
$ws.Range('A13').Value = '陳逸(c_personid=14542)的入仕方式是？'
$ws.Range('B13').Value = '恩蔭'
$ws.Range('A14').Value = '劉伸(c_personid=557888的籍貫是否為陽州？'
$ws.Range('B14').Value = '否'
$ws.Range('A15').Value = '陳邦光(c_personid=74)的入仕方式是？'
$ws.Range('B15').Value = '進士'
$ws.Range('A16').Value = '單瑞龍(c_personid=368518)的籍貫是？'
$ws.Range('B16').Value = '錢塘'
$ws.Range('A17').Value = '陸壑(c_personid=48095)的入仕方式是？'
$ws.Range('B17').Value = '進士'
$ws.Range('A18').Value = '林叢桂(c_personid=624091)的入仕方式是？'
$ws.Range('B18').Value = '舉人'
$ws.Range('A19').Value = '李有普(c_personid=620769)的入仕方式是？'
$ws.Range('B19').Value = '恩蔭'
$ws.Range('A20').Value = '王文烱(c_personid=327293的籍貫是否為平州？'
$ws.Range('B20').Value = '否'
$ws.Range('A21').Value = '舒嘉猷(c_personid=90244的籍貫是否為慶州？'
$ws.Range('B21').Value = '否'
$ws.Range('A22').Value = '徐淮(c_personid=385952)的入仕方式是？'
$ws.Range('B22').Value = '進士'
$ws.Range('A23').Value = '張公邵(c_personid=26542)的入仕方式是？'
$ws.Range('B23').Value = '恩蔭'
$ws.Range('A24').Value = '周春(c_personid=33209的籍貫是否為寧州？'
$ws.Range('B24').Value = '否'
$ws.Range('A25').Value = '王杲(c_personid=39477)的籍貫是？'
$ws.Range('B25').Value = '齊州'
$ws.Range('A26').Value = '胡永榮(c_personid=648279)的入仕方式是？'
$ws.Range('B26').Value = '舉人'
$ws.Range('A27').Value = '呂弼康(c_personid=43017)的籍貫是？'
$ws.Range('B27').Value = '餘干'
$ws.Range('A28').Value = '易俊(c_personid=59256的籍貫是否為國州？'
$ws.Range('B28').Value = '否'
$ws.Range('A29').Value = '徐子端(c_personid=27782)的入仕方式是？'
$ws.Range('B29').Value = '恩蔭'
$ws.Range('A30').Value = '包景寧(c_personid=589572)的入仕方式是？'
$ws.Range('B30').Value = '舉人'
$ws.Range('A31').Value = '劉運隆(c_personid=369916)的入仕方式是？'
$ws.Range('B31').Value = '進士'
$ws.Range('A32').Value = '王鵬運(c_personid=54979)的入仕方式是？'
$ws.Range('B32').Value = '舉人'
$ws.Range('A33').Value = '胡鶴(c_personid=300994)的籍貫是？'
$ws.Range('B33').Value = '歙縣'
$ws.Range('A34').Value = '陳徽言(c_personid=82250的籍貫是否為川州？'
$ws.Range('B34').Value = '否'
$ws.Range('A35').Value = '孔羲仲(c_personid=126186)的入仕方式是？'
$ws.Range('B35').Value = '恩蔭'
$ws.Range('A36').Value = '陳宗鳳(c_personid=670003)的入仕方式是？'
$ws.Range('B36').Value = '舉人'
$ws.Range('A37').Value = '某徵(c_personid=555494)的入仕方式是？'
$ws.Range('B37').Value = '進士'
$ws.Range('A38').Value = '劉廷聞(c_personid=586841)的入仕方式是？'
$ws.Range('B38').Value = '舉人'
$ws.Range('A39').Value = '張淑躬(c_personid=608044)的入仕方式是？'
$ws.Range('B39').Value = '舉人'
$ws.Range('A40').Value = '李名瑚(c_personid=618707)的入仕方式是？'
$ws.Range('B40').Value = '恩蔭'
$ws.Range('A41').Value = '尹光楨(c_personid=602462)的籍貫是？'
$ws.Range('B41').Value = '新城'
$ws.Range('A42').Value = '程再伊(c_personid=559998)的籍貫是？'
$ws.Range('B42').Value = '鄱陽'
$ws.Range('A43').Value = '吳公瑾(c_personid=593506)的入仕方式是？'
$ws.Range('B43').Value = '恩蔭'
$ws.Range('A44').Value = '辛本婺(c_personid=661823)的籍貫是？'
$ws.Range('B44').Value = '蓬萊'
$ws.Range('A45').Value = '蔡雲吉(c_personid=493073)的入仕方式是？'
$ws.Range('B45').Value = '舉人'
$ws.Range('A46').Value = '胡實(c_personid=11871)的入仕方式是？'
$ws.Range('B46').Value = '恩蔭'
$ws.Range('A47').Value = '張曾垿(c_personid=347008)的入仕方式是？'
$ws.Range('B47').Value = '進士'
$ws.Range('A48').Value = '吳秉翰(c_personid=592604的籍貫是否為平府直轄地方？'
$ws.Range('B48').Value = '否'
$ws.Range('A49').Value = '陳成務(c_personid=541971)的籍貫是？'
$ws.Range('B49').Value = '晉江'
$ws.Range('A50').Value = '石萬寶(c_personid=643103)的籍貫是？'
$ws.Range('B50').Value = '如皋'
$ws.Range('A51').Value = '松壽(c_personid=366024)的入仕方式是？'
$ws.Range('B51').Value = '進士'
